# Auto-generated update of cryptos price/volume cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.231.63"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.690.67"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.18"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.48%  "
$ws.Range("E9").Value = "  +4.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0629"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.928.26"
$ws.Range("D13").Value = "1.694.23"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "27.226.07"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.25"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "1.554.81"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "1.836.83"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.67"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "0.0₆0112"
$ws.Range("E48").Value = "  +7.47%  "
$ws.Range("E49").Value = "  +6.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.40"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.53%  "
$ws.Range("E51").Value = "  +1.18%  "

Write-Output "done"